# Append three new paragraphs after the last paragraph of the document body
# ("We can append a list to another list by using + or extend method"):
#   1. An empty paragraph (body-text size, 12pt / sz=24)
#   2. A bold, 14pt (sz=28) paragraph containing "07.03."
#   3. An empty bold 14pt (sz=28) paragraph
#
# Built/inserted as raw WordprocessingML via Range.InsertXML so the emitted
# paragraphs contain no stray empty runs, matching how Word itself authors
# blank paragraphs (pPr/rPr only, no <w:r>).

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$blankBody = '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
        '<w:spacing w:after="0"/>' +
        '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
    '</w:pPr>' +
'</w:p>'

$dateHeading = '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
        '<w:spacing w:after="0"/>' +
        '<w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
        '<w:t>07.03.</w:t>' +
    '</w:r>' +
'</w:p>'

$blankHeading = '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
        '<w:spacing w:after="0"/>' +
        '<w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
    '</w:pPr>' +
'</w:p>'

$newParagraphsXml = $blankBody + $dateHeading + $blankHeading

$insertionPoint = $d.Content
$insertionPoint.Collapse(0)
$insertionPoint.InsertXML($newParagraphsXml)
